$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E12" = 18.16350000000002
    "E27" = 16.68279999999999
    "E32" = 16.70239999999998
    "E36" = 17.23460000000001
    "E38" = 16.72329999999999
    "E46" = 17.0975
    "E54" = 16.6382
    "E55" = 16.64280000000001
    "E56" = 16.7377
    "E67" = 17.15940000000002
    "E69" = 17.17860000000003
    "E72" = 16.5259
    "E83" = 16.703
    "E86" = 16.68580000000001
    "E91" = 18.48800000000003
    "E93" = 17.36780000000002
    "E99" = 16.5099
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
